$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data refresh
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.318.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.88%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.427.79'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.98%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.01'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.576'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.13%  '

$ws.Range("E9").Value = '  +4.77%  '

$ws.Range("E10").Value = '  +3.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.359'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.45%  '

$ws.Range("E12").Value = '  -2.31%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.92'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.860.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.03%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.242.92'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.83%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000139'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.75%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.426.30'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.09%  '

$ws.Range("E18").Value = '  +5.96%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '335.46'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.08%  '

$ws.Range("E22").Value = '  -0.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.49%  '

$ws.Range("E24").Value = '  +3.37%  '

$ws.Range("E25").Value = '  +1.36%  '

$ws.Range("E26").Value = '  +0.08%  '

$ws.Range("E27").Value = '  -0.31%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0787'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.45%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.79'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.30'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '169.37'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.82%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.78'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.67%  '

$ws.Range("E33").Value = '  +1.31%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.31'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.23'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.80%  '

$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.62'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '39.84'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '322.29'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.90%  '

$ws.Range("E41").Value = '  +10.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '142.70'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0526'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.75%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0960'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.97%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.64'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.411'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.15%  '

$ws.Range("E49").Value = '  +1.60%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.90'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.83%  '

$ws.Range("E51").Value = '  -0.07%  '
